$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "is_Mandatory"
$ws.Range("E1").Value = "default_value"
$ws.Range("F1").Value = "enable_check_for_invalid"

$ws.Range("D2").Select()
